$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H116").Value = 1887.8125
$ws.Range("I116").Value = 1920.5
$ws.Range("J116").Value = 1833.3334
$ws.Range("K116").Value = 1920.5
$ws.Range("L116").Value = 1833.3334
$ws.Range("M116").Value = 1521.5
$ws.Range("N116").Value = -8717.3334
$ws.Range("H129").Value = 563903.9
$ws.Range("I129").Value = 636.7778
$ws.Range("K129").Value = 1910.3334
$ws.Range("M129").Value = 3089.6666
$ws.Range("H137").Value = 1258.1455
$ws.Range("I137").Value = 1064.5581
$ws.Range("J137").Value = 1951.8334
$ws.Range("K137").Value = 3193.6743
$ws.Range("L137").Value = 5855.5002
$ws.Range("M137").Value = -643.6742999999997
$ws.Range("N137").Value = -10955.5002
$ws.Range("H138").Value = 1846.9667
$ws.Range("I138").Value = 1505.3667
$ws.Range("J138").Value = 2188.5667
$ws.Range("K138").Value = 4516.1001
$ws.Range("L138").Value = 6565.7001
$ws.Range("M138").Value = 623.8999000000003
$ws.Range("N138").Value = -16845.7001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1273.32
$ws.Range("I2").Value = 920.1667
$ws.Range("J2").Value = 1599.3077
$ws.Range("K2").Value = 920.1667
$ws.Range("L2").Value = 1599.3077
$ws.Range("M2").Value = -807.1667
$ws.Range("N2").Value = -1825.3077
$ws.Range("H5").Value = 92.75
$ws.Range("I5").Value = 90.5
$ws.Range("J5").Value = 95
$ws.Range("K5").Value = 90.5
$ws.Range("L5").Value = 95
$ws.Range("M5").Value = 21.5
$ws.Range("N5").Value = -319
$ws.Range("H32").Value = 1196.44
$ws.Range("I32").Value = 1119.4946
$ws.Range("J32").Value = 2218.7144
$ws.Range("K32").Value = 1119.4946
$ws.Range("L32").Value = 2218.7144
$ws.Range("M32").Value = -832.4946
$ws.Range("N32").Value = -2792.7144
$ws.Range("H74").Value = 703.18604
$ws.Range("I74").Value = 585.925
$ws.Range("J74").Value = 2266.6667
$ws.Range("K74").Value = 585.925
$ws.Range("L74").Value = 2266.6667
$ws.Range("M74").Value = 288.075
$ws.Range("N74").Value = -4014.6667
$ws.Range("H77").Value = 703.18604
$ws.Range("I77").Value = 585.925
$ws.Range("J77").Value = 2266.6667
$ws.Range("K77").Value = 2929.625
$ws.Range("L77").Value = 11333.3335
$ws.Range("M77").Value = 1438.375
$ws.Range("N77").Value = -20069.3335
$ws.Range("H116").Value = 1273.32
$ws.Range("I116").Value = 920.1667
$ws.Range("J116").Value = 1599.3077
$ws.Range("K116").Value = 920.1667
$ws.Range("L116").Value = 1599.3077
$ws.Range("M116").Value = 1373.8333
$ws.Range("N116").Value = -6187.3077
$ws.Range("H132").Value = 3485.15
$ws.Range("I132").Value = 3749.3865
$ws.Range("J132").Value = 2758.5
$ws.Range("K132").Value = 11248.1595
$ws.Range("L132").Value = 8275.5
$ws.Range("M132").Value = -8718.1595
$ws.Range("N132").Value = -13335.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1273.32
$ws.Range("I3").Value = 920.1667
$ws.Range("J3").Value = 1599.3077
$ws.Range("K3").Value = 920.1667
$ws.Range("L3").Value = 1599.3077
$ws.Range("M3").Value = -806.1667
$ws.Range("N3").Value = -1827.3077
$ws.Range("H4").Value = 92.75
$ws.Range("I4").Value = 90.5
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 90.5
$ws.Range("L4").Value = 95
$ws.Range("M4").Value = 24.5
$ws.Range("N4").Value = -325
$ws.Range("H20").Value = 2215
$ws.Range("I20").Value = 2027.125
$ws.Range("K20").Value = 2027.125
$ws.Range("M20").Value = -1780.125
$ws.Range("H99").Value = 1073.3334
$ws.Range("I99").Value = 1095
$ws.Range("J99").Value = 900
$ws.Range("K99").Value = 1095
$ws.Range("L99").Value = 900
$ws.Range("M99").Value = 403
$ws.Range("N99").Value = -3896
$ws.Range("H105").Value = 3297.2727
$ws.Range("I105").Value = 2490
$ws.Range("J105").Value = 6930
$ws.Range("K105").Value = 2490
$ws.Range("L105").Value = 6930
$ws.Range("M105").Value = -743
$ws.Range("N105").Value = -10424

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 51600
$ws.Range("I141").Value = 29000
$ws.Range("J141").Value = 57250
$ws.Range("K141").Value = 29000
$ws.Range("L141").Value = 57250
$ws.Range("M141").Value = -23820
$ws.Range("N141").Value = -67610

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I34").Value = 2000
$ws.Range("K34").Value = 6000
$ws.Range("M34").Value = -5916
$ws.Range("H36").Value = 1063.2307
$ws.Range("I36").Value = 1063.2307
$ws.Range("K36").Value = 3189.6921
$ws.Range("M36").Value = -3020.6921
$ws.Range("H86").Value = 606.6667
$ws.Range("I86").Value = 820
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 2460
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -1274
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 606.6667
$ws.Range("I89").Value = 820
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 7380
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = -1452
$ws.Range("N89").Value = -16356
$ws.Range("H98").Value = 948.5
$ws.Range("I98").Value = 600
$ws.Range("J98").Value = 1157.6
$ws.Range("K98").Value = 1800
$ws.Range("L98").Value = 3472.8
$ws.Range("M98").Value = -302
$ws.Range("N98").Value = -6468.799999999999
$ws.Range("H131").Value = 1427261.6
$ws.Range("J131").Value = 1838847.5
$ws.Range("L131").Value = 5516542.5
$ws.Range("N131").Value = -5526622.5
$ws.Range("H132").Value = 2126.375
$ws.Range("I132").Value = 1241.625
$ws.Range("J132").Value = 2568.75
$ws.Range("K132").Value = 11174.625
$ws.Range("L132").Value = 23118.75
$ws.Range("M132").Value = -8644.625
$ws.Range("N132").Value = -28178.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 27779744
$ws.Range("I113").Value = 31252038
$ws.Range("K113").Value = 31252038
$ws.Range("M113").Value = -31249868
$ws.Range("H131").Value = 22000
$ws.Range("J131").Value = 22000
$ws.Range("L131").Value = 22000
$ws.Range("N131").Value = -32080

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1735
$ws.Range("I16").Value = 554.2857
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 554.2857
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -384.2857
$ws.Range("N16").Value = -10340
$ws.Range("H55").Value = 384.72726
$ws.Range("I55").Value = 449.75
$ws.Range("J55").Value = 347.57144
$ws.Range("K55").Value = 449.75
$ws.Range("L55").Value = 347.57144
$ws.Range("M55").Value = -276.75
$ws.Range("N55").Value = -693.5714399999999
$ws.Range("H101").Value = 20344.6
$ws.Range("J101").Value = 20344.6
$ws.Range("L101").Value = 20344.6
$ws.Range("N101").Value = -26834.6
$ws.Range("H102").Value = 50000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -56490
$ws.Range("H105").Value = 24815
$ws.Range("J105").Value = 24815
$ws.Range("L105").Value = 24815
$ws.Range("N105").Value = -31803
$ws.Range("H132").Value = 7019.2856
$ws.Range("I132").Value = 10053.909
$ws.Range("K132").Value = 30161.727
$ws.Range("M132").Value = -27631.727
$ws.Range("H136").Value = 2465.3438
$ws.Range("I136").Value = 2424.9443
$ws.Range("J136").Value = 2683.5
$ws.Range("K136").Value = 7274.8329
$ws.Range("L136").Value = 8050.5
$ws.Range("M136").Value = -4724.8329
$ws.Range("N136").Value = -13150.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 888.9655
$ws.Range("I132").Value = 846.17645
$ws.Range("J132").Value = 1200.7142
$ws.Range("K132").Value = 2538.52935
$ws.Range("L132").Value = 3602.1426
$ws.Range("M132").Value = -8.52935000000025
$ws.Range("N132").Value = -8662.142599999999
